# Append a new data row (row 29) to the active sheet, mirroring the
# latest Adafruit IO reading (same values as the preceding row, 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
# Leading apostrophe forces Excel to store this numeric-looking reading
# as text (matching the existing text-typed "Value" column) rather than
# silently coercing it into a number.
$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
